$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new data rows (id/name pairs) below the existing table data.
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "c1"

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "c2"

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "c3"

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "c4"

# Grow Table1 so it covers the newly added rows.
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:B7"))

# Move/restore the active cell selection to match the saved view state.
[void]$ws.Range("B13").Select()
